$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "username"/"password" row (row 2), shifting subsequent rows up
$ws.Rows.Item(2).Delete()

# Update the selection to reflect the new active cell
$ws.Range("A2").Select()
